$d = $word.ActiveDocument
$CR = [char]13

# -----------------------------------------------------------------
# Part 1: add a new "{#render}" paragraph (styled "Normale") right
# before the bold "Trimestre" heading paragraph.
# -----------------------------------------------------------------
$trimestreIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq ("Trimestre" + $CR)) {
        $trimestreIdx = $i
        break
    }
}

$anchor1 = $d.Paragraphs.Item($trimestreIdx - 1)
$anchor1.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($trimestreIdx)
$newPara.Range.Text = "{#render}"
$newPara.Style = "Normale"

# Touch Bold on/off to materialise an explicit (empty) run-properties
# element on the new run, matching the target markup's <w:rPr/>.
$renderFind = $d.Content
$renderFind.Find.Execute("{#render}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$renderFind.Font.Bold = $true
$renderFind.Font.Bold = $false

# -----------------------------------------------------------------
# Part 2: fill in the blank paragraph that follows "{/argomenti_q2}"
# with "{/render}" (split into "{", "/render", "}" runs) and reset
# its left indent from 720 to 0.
# -----------------------------------------------------------------
$closeIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq ("{/argomenti_q2}" + $CR)) {
        $closeIdx = $i
        break
    }
}

$blankPara = $d.Paragraphs.Item($closeIdx + 1)
$blankPara.LeftIndent = 0
$blankPara.Range.Text = "{/render}"

$fullRun = $d.Paragraphs.Item($closeIdx + 1).Range
$fullRun.Font.Name = "Arial"
$fullRun.Font.NameFarEast = "Arial"
$fullRun.Font.NameBi = "Arial"
$fullRun.Font.Size = 12
$fullRun.Font.SizeBi = 12

# Re-apply the same formatting to the "/render" substring only; this
# forces Word to split the single run into three runs - "{", "/render"
# and "}" - while preserving identical run formatting on each.
$splitFind = $d.Content
$splitFind.Find.Execute("/render", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitFind.Font.Name = "Arial"
$splitFind.Font.NameFarEast = "Arial"
$splitFind.Font.NameBi = "Arial"
$splitFind.Font.Size = 12
$splitFind.Font.SizeBi = 12

Write-Output "edit complete"
